$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-66 down to 25-67
$ws.Rows.Item(24).Insert()

# Copy the (now shifted) original row 25 formatting/values into the new blank row 24,
# so the static columns (A, B, C, E-L, Q, R, T) match the rest of the dataset.
$ws.Range("A25:T25").Copy()
$ws.Range("A24").PasteSpecial()

# Now overwrite this week's new observation values in row 24
$ws.Range("D24").Value = 44708
$ws.Range("M24").Value = 20
$ws.Range("N24").Value = 30000
$ws.Range("O24").Value = 30000
$ws.Range("P24").Value = 30000
$ws.Range("S24").Value = 1500

$excel.CutCopyMode = $false
